{"js": "const pairs = [\n  [\"78\u00d745=\", \"31\u00d738=\"],\n  [\"12\u00d739=\", \"50\u00d784=\"],\n  [\"39\u00d776=\", \"30\u00d762=\"],\n  [\"71\u00d765=\", \"45\u00d792=\"],\n  [\"96\u00d719=\", \"17\u00d756=\"],\n  [\"88\u00d798=\", \"93\u00d720=\"],\n  [\"27\u00d780=\", \"22\u00d720=\"],\n  [\"95\u00d734=\", \"62\u00d741=\"],\n  [\"20\u00d741=\", \"28\u00d767=\"],\n  [\"56\u00d779=\", \"57\u00d732=\"],\n  [\"15\u00d732=\", \"62\u00d796=\"],\n  [\"53\u00d798=\", \"81\u00d757=\"],\n  [\"43\u00d789=\", \"99\u00d764=\"],\n  [\"79\u00d725=\", \"38\u00d712=\"],\n  [\"66\u00d737=\", \"77\u00d727=\"],\n  [\"19\u00d773=\", \"67\u00d784=\"],\n  [\"20\u00d713=\", \"99\u00d792=\"],\n  [\"27\u00d779=\", \"75\u00d737=\"],\n  [\"80\u00d782=\", \"40\u00d721=\"],\n  [\"42\u00d790=\", \"82\u00d724=\"],\n  [\"89\u00d771=\", \"71\u00d782=\"],\n  [\"52\u00d769=\", \"93\u00d723=\"],\n  [\"53\u00d767=\", \"91\u00d780=\"],\n  [\"39\u00d737=\", \"20\u00d717=\"],\n  [\"57\u00d799=\", \"14\u00d744=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"78\u00d745=\", \"31\u00d738=\"),\n    @(\"12\u00d739=\", \"50\u00d784=\"),\n    @(\"39\u00d776=\", \"30\u00d762=\"),\n    @(\"71\u00d765=\", \"45\u00d792=\"),\n    @(\"96\u00d719=\", \"17\u00d756=\"),\n    @(\"88\u00d798=\", \"93\u00d720=\"),\n    @(\"27\u00d780=\", \"22\u00d720=\"),\n    @(\"95\u00d734=\", \"62\u00d741=\"),\n    @(\"20\u00d741=\", \"28\u00d767=\"),\n    @(\"56\u00d779=\", \"57\u00d732=\"),\n    @(\"15\u00d732=\", \"62\u00d796=\"),\n    @(\"53\u00d798=\", \"81\u00d757=\"),\n    @(\"43\u00d789=\", \"99\u00d764=\"),\n    @(\"79\u00d725=\", \"38\u00d712=\"),\n    @(\"66\u00d737=\", \"77\u00d727=\"),\n    @(\"19\u00d773=\", \"67\u00d784=\"),\n    @(\"20\u00d713=\", \"99\u00d792=\"),\n    @(\"27\u00d779=\", \"75\u00d737=\"),\n    @(\"80\u00d782=\", \"40\u00d721=\"),\n    @(\"42\u00d790=\", \"82\u00d724=\"),\n    @(\"89\u00d771=\", \"71\u00d782=\"),\n    @(\"52\u00d769=\", \"93\u00d723=\"),\n    @(\"53\u00d767=\", \"91\u00d780=\"),\n    @(\"39\u00d737=\", \"20\u00d717=\"),\n    @(\"57\u00d799=\", \"14\u00d744=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"Replacement failed for $oldText\"\n    }\n}"}
